$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) and volume-change (column E) updates
# Values are kept as text (matching the source workbook, which stores
# these columns as text strings) by forcing the cell number format to
# "@" (Text) before assigning the value, avoiding Excel's automatic
# numeric/percentage conversion.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "286.84"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.27%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.33"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.10%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.917"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.46%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.50%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.247"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.30%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.368"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "15.14%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9135"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.24%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1570"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.82%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06714"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "31.31%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07698"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.74%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02982"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.35%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08978"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.02%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001605"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.52%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006545"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.32%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006025"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.57%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.470"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.06%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.396"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.69%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.45%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.69%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.965"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.44%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.37%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1520"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "10.12%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001185"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.55%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004334"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "12.54%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001179"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.77%"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-15.73%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04169"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.93%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006771"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.77%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1414"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "20.58%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002159"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.36%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01241"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.08%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005569"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.04%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-7.03%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.25%"
